# Add a blank "comment" row above the header row on every sheet template.
# The header row (previously row 3) moves up to row 2, leaving row 1 blank
# for a future comment, and any data rows stay where they are.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $headerRow = 3
    $newHeaderRow = 2

    $oldRange = "A" + $headerRow + ":X" + $headerRow
    $newRange = "A" + $newHeaderRow

    # Move the header row's values/styles up one row (row 3 -> row 2),
    # without disturbing any rows below the header (e.g. sample data on
    # row 4 of the Individuals sheet must stay put).
    $ws.Range($oldRange).Cut($ws.Range($newRange))

    # The old header row is now a stray, fully-blank row (no values, no
    # left-over cell styling) sitting between the header and any data
    # below it - clear it out completely so it disappears from the
    # saved sheetData entirely.
    $oldRowRange = $ws.Rows($headerRow.ToString() + ":" + $headerRow.ToString())
    $oldRowRange.Delete()
    $oldRowRange.Insert()
    $oldRowRange.Clear()
}

# Select the new blank first row on every sheet (entire row, A:XFD) - this
# is where a comment will be typed in - and make "Groups" the active tab.
$groups = $wb.Worksheets.Item("Groups")
$individuals = $wb.Worksheets.Item("Individuals")

$individuals.Activate()
$individuals.Range("A2:XFD2").Select()

$groups.Activate()
$groups.Range("A1:XFD1").Select()
